# Update columns E and F for a set of rows on the active worksheet.
# These updates bring the "current" (E) and "projected" (F) totals in
# line with the latest refreshed "previous" (D) totals following the
# automated `make clean download import process` run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    6   = 242342.31
    9   = 124860.15
    16  = 67825.58
    32  = 177069.06
    50  = 27381.49
    64  = 13862.58
    79  = 164884.39
    84  = 52266.74
    86  = 615686.6
    95  = 95236
    100 = 58399.56
    105 = 68499.13
    113 = 191189.79
    114 = 21036.5
    142 = 76910.71000000001
    143 = 134325.58
    157 = 21685.48
    166 = 2294.65
    181 = 175063
}

foreach ($row in $updates.Keys) {
    $value = $updates[$row]
    $ws.Range("E$row").Value = $value
    $ws.Range("F$row").Value = $value
}
